# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to the "想去人数" (want-to-go count) column F
# and marks one event as "不可售" (not for sale) in column G
# across all affected worksheets.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value  = 7669
$ws1.Range("F3").Value  = 7669
$ws1.Range("F5").Value  = 7847
$ws1.Range("F8").Value  = 31
$ws1.Range("F9").Value  = 6612
$ws1.Range("F10").Value = 3367
$ws1.Range("F12").Value = 3712
$ws1.Range("F13").Value = 43
$ws1.Range("F14").Value = 44
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 464
$ws1.Range("F20").Value = 28
$ws1.Range("F22").Value = 326
$ws1.Range("F23").Value = 3826
$ws1.Range("F25").Value = 369
$ws1.Range("F26").Value = 953
$ws1.Range("F28").Value = 1467
$ws1.Range("F30").Value = 54
$ws1.Range("F31").Value = 2742
$ws1.Range("F32").Value = 1803
$ws1.Range("F33").Value = 31
$ws1.Range("F36").Value = 3652
$ws1.Range("F37").Value = 307
$ws1.Range("F38").Value = 279
$ws1.Range("F41").Value = 534
$ws1.Range("G41").Value = "不可售"
$ws1.Range("F42").Value = 1412
$ws1.Range("F45").Value = 636
$ws1.Range("F46").Value = 1

# ------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F13").Value = 89
$ws2.Range("F17").Value = 24

# ------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 134

# ------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value  = 134
$ws4.Range("F5").Value  = 7669
$ws4.Range("F6").Value  = 7669
$ws4.Range("F7").Value  = 7847
$ws4.Range("F9").Value  = 31
$ws4.Range("F10").Value = 6612
$ws4.Range("F11").Value = 3367
$ws4.Range("F12").Value = 3712
$ws4.Range("F13").Value = 44
$ws4.Range("F16").Value = 58
$ws4.Range("F17").Value = 464
$ws4.Range("F19").Value = 28
$ws4.Range("F22").Value = 326
$ws4.Range("F23").Value = 3826
$ws4.Range("F27").Value = 369
$ws4.Range("F28").Value = 953
$ws4.Range("F30").Value = 1467
$ws4.Range("F32").Value = 54
$ws4.Range("F33").Value = 2742
$ws4.Range("F34").Value = 1803
$ws4.Range("F35").Value = 31
$ws4.Range("F38").Value = 3652
$ws4.Range("F39").Value = 307
$ws4.Range("F40").Value = 279
$ws4.Range("F44").Value = 534
$ws4.Range("G44").Value = "不可售"
$ws4.Range("F45").Value = 24
$ws4.Range("F46").Value = 1412
$ws4.Range("F50").Value = 636
